$d = $word.ActiveDocument

$d.Content.Find.Execute("838×6=", $true, $false, $false, $false, $false, $true, 1, $false, "927×6=", 2) | Out-Null
$d.Content.Find.Execute("175×4=", $true, $false, $false, $false, $false, $true, 1, $false, "278×8=", 2) | Out-Null
$d.Content.Find.Execute("531×2=", $true, $false, $false, $false, $false, $true, 1, $false, "900×2=", 2) | Out-Null
$d.Content.Find.Execute("153×5=", $true, $false, $false, $false, $false, $true, 1, $false, "550×5=", 2) | Out-Null
$d.Content.Find.Execute("728×7=", $true, $false, $false, $false, $false, $true, 1, $false, "598×2=", 2) | Out-Null
$d.Content.Find.Execute("521×9=", $true, $false, $false, $false, $false, $true, 1, $false, "658×2=", 2) | Out-Null
$d.Content.Find.Execute("804×6=", $true, $false, $false, $false, $false, $true, 1, $false, "991×9=", 2) | Out-Null
$d.Content.Find.Execute("557×4=", $true, $false, $false, $false, $false, $true, 1, $false, "580×2=", 2) | Out-Null
$d.Content.Find.Execute("965×7=", $true, $false, $false, $false, $false, $true, 1, $false, "641×6=", 2) | Out-Null
$d.Content.Find.Execute("470×9=", $true, $false, $false, $false, $false, $true, 1, $false, "669×9=", 2) | Out-Null
$d.Content.Find.Execute("365×9=", $true, $false, $false, $false, $false, $true, 1, $false, "302×2=", 2) | Out-Null
$d.Content.Find.Execute("150×7=", $true, $false, $false, $false, $false, $true, 1, $false, "935×7=", 2) | Out-Null
$d.Content.Find.Execute("559×2=", $true, $false, $false, $false, $false, $true, 1, $false, "933×4=", 2) | Out-Null
$d.Content.Find.Execute("876×6=", $true, $false, $false, $false, $false, $true, 1, $false, "876×4=", 2) | Out-Null
$d.Content.Find.Execute("340×2=", $true, $false, $false, $false, $false, $true, 1, $false, "618×9=", 2) | Out-Null
$d.Content.Find.Execute("376×8=", $true, $false, $false, $false, $false, $true, 1, $false, "135×4=", 2) | Out-Null
$d.Content.Find.Execute("846×5=", $true, $false, $false, $false, $false, $true, 1, $false, "421×4=", 2) | Out-Null
$d.Content.Find.Execute("861×3=", $true, $false, $false, $false, $false, $true, 1, $false, "443×6=", 2) | Out-Null
$d.Content.Find.Execute("227×9=", $true, $false, $false, $false, $false, $true, 1, $false, "154×9=", 2) | Out-Null
$d.Content.Find.Execute("358×5=", $true, $false, $false, $false, $false, $true, 1, $false, "654×7=", 2) | Out-Null
$d.Content.Find.Execute("559×6=", $true, $false, $false, $false, $false, $true, 1, $false, "369×8=", 2) | Out-Null
$d.Content.Find.Execute("342×3=", $true, $false, $false, $false, $false, $true, 1, $false, "364×5=", 2) | Out-Null
$d.Content.Find.Execute("855×7=", $true, $false, $false, $false, $false, $true, 1, $false, "220×6=", 2) | Out-Null
$d.Content.Find.Execute("233×4=", $true, $false, $false, $false, $false, $true, 1, $false, "502×7=", 2) | Out-Null
$d.Content.Find.Execute("159×9=", $true, $false, $false, $false, $false, $true, 1, $false, "399×7=", 2) | Out-Null
